$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Daily level" sheet: rename the Hydration header, rewrite row 2, and add
# 9 new rows (3-11) of tracked data plus a "Total Calories" formula column.
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily level")

# Header: "Hydration" -> "Hydration(liters)"
$daily.Range("C1").Value = "Hydration(liters)"

# Row 2 updates (hydration now fractional liters; fruit/veg + new calories total)
$daily.Range("C2").Value2 = 2.5
$daily.Range("E2").Value2 = 350
$daily.Range("F2").Value2 = 218
$daily.Range("G2").Formula = "=+E2+F2"

# New data rows 3-11
$stepCounts = @(7828, 6000, 7500, 6500, 6472, 7200, 6738, 6389, 7863)
$hydration  = @(2, 3, 2.5, 3, 2.8, 2.6, 2, 3, 3)
$sleep      = @(9, 8, 7, 7, 6.5, 7, 8, 7.5, 8)
$fruit      = @(320, 420, 360, 290, 209, 467, 527, 268, 638)
$veg        = @(250, 290, 220, 350, 420, 299, 201, 479, 493)
$dates      = @(44043, 44047, 44058, 44063, 44066, 44066, 44063, 44044, 44064)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 3
    $daily.Cells.Item($row, 1).Value2 = $dates[$i]
    $daily.Cells.Item($row, 2).Value2 = $stepCounts[$i]
    $daily.Cells.Item($row, 3).Value2 = $hydration[$i]
    $daily.Cells.Item($row, 4).Value2 = $sleep[$i]
    $daily.Cells.Item($row, 5).Value2 = $fruit[$i]
    $daily.Cells.Item($row, 6).Value2 = $veg[$i]
}

# Match date number formatting (style) used by the existing A2 cell
$daily.Range("A2").Copy()
$daily.Range("A3:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Total Calories formula column for the new rows
$daily.Range("G3:G11").Formula = "=+E3+F3"

# ---------------------------------------------------------------------------
# "Group Member" sheet: add 10 rows of invite/accept/reject tracking data.
# ---------------------------------------------------------------------------
$groupMember = $wb.Worksheets.Item("Group Member")

# Columns: A=Date Invited, B=Date Accepted Invite, C=Date Rejected Invite, D=Date Left Group
# Written column-by-column (matching the original author's entry order, which
# determines the shared-string table order: Pending, Rejected, pending, Accepted).
$colA = @(44066, 44066, 44063, 44044, 44064, 44067, 44043, 44047, 44058, 44063)
for ($i = 0; $i -lt $colA.Length; $i++) {
    $groupMember.Cells.Item($i + 2, 1).Value2 = $colA[$i]
}

$groupMember.Cells.Item(2, 2).Value = "Pending"
$groupMember.Cells.Item(3, 2).Value = "Pending"
$groupMember.Cells.Item(4, 2).Value2 = 44064
$groupMember.Cells.Item(5, 2).Value2 = 44044
$groupMember.Cells.Item(6, 2).Value2 = 44065
$groupMember.Cells.Item(7, 2).Value = "Pending"
$groupMember.Cells.Item(8, 2).Value2 = 44044
$groupMember.Cells.Item(9, 2).Value2 = 44048
$groupMember.Cells.Item(10, 2).Value2 = 44063
$groupMember.Cells.Item(11, 2).Value = "Rejected"

$groupMember.Cells.Item(2, 3).Value = "pending"
$groupMember.Cells.Item(3, 3).Value = "Pending"
$groupMember.Cells.Item(4, 3).Value = "Accepted"
$groupMember.Cells.Item(5, 3).Value = "Accepted"
$groupMember.Cells.Item(6, 3).Value = "Accepted"
$groupMember.Cells.Item(7, 3).Value = "pending"
$groupMember.Cells.Item(8, 3).Value = "Accepted"
$groupMember.Cells.Item(9, 3).Value = "Accepted"
$groupMember.Cells.Item(10, 3).Value = "Accepted"
$groupMember.Cells.Item(11, 3).Value2 = 44064

$groupMember.Cells.Item(2, 4).Value = "pending"
$groupMember.Cells.Item(3, 4).Value = "pending"
$groupMember.Cells.Item(4, 4).Value = "Accepted"
$groupMember.Cells.Item(5, 4).Value = "Accepted"
$groupMember.Cells.Item(6, 4).Value = "Accepted"
$groupMember.Cells.Item(7, 4).Value = "pending"
$groupMember.Cells.Item(8, 4).Value2 = 44063
$groupMember.Cells.Item(9, 4).Value = "Accepted"
$groupMember.Cells.Item(10, 4).Value = "Accepted"
$groupMember.Cells.Item(11, 4).Value = "Accepted"

# Apply the workbook's existing short-date style (the one already used by
# "Daily level"!A2) to every date-valued cell on the new rows.
$daily.Range("A2").Copy()
$groupMember.Range("A2:A11").PasteSpecial(-4122)
$groupMember.Range("B4:B6").PasteSpecial(-4122)
$groupMember.Range("B8:B10").PasteSpecial(-4122)
$groupMember.Range("C11").PasteSpecial(-4122)
$groupMember.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$groupMember.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping: "Daily level" becomes the active tab,
# matching the author's final view state when they saved the workbook.
# ---------------------------------------------------------------------------
$groupMember.Range("E11").Select()

$daily.Activate()
$daily.Range("C17").Select()
